$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.287956476211548
$ws.Range("B1").Value = 2.937998533248901
$ws.Range("C1").Value = 5.273559093475342
$ws.Range("D1").Value = 1.847816824913025
$ws.Range("E1").Value = 1.0152188539505
